$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 99999.8
$ws.Range("J3").Value = 99999.8
$ws.Range("L3").Value = 99999.8
$ws.Range("N3").Value = -100227.8
$ws.Range("H40").Value = 4482.3335
$ws.Range("I40").Value = 4482.3335
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4482.3335
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4307.3335
$ws.Range("N40").ClearContents()
$ws.Range("H95").Value = 56500
$ws.Range("J95").Value = 56500
$ws.Range("L95").Value = 56500
$ws.Range("N95").Value = -61992
$ws.Range("H97").Value = 4049.25
$ws.Range("I97").Value = 1500
$ws.Range("J97").Value = 4899
$ws.Range("K97").Value = 4500
$ws.Range("L97").Value = 14697
$ws.Range("M97").Value = -4004
$ws.Range("N97").Value = -15689
$ws.Range("H99").Value = 3018.1667
$ws.Range("I99").Value = 234.42857
$ws.Range("J99").Value = 6915.4
$ws.Range("K99").Value = 703.28571
$ws.Range("L99").Value = 20746.2
$ws.Range("M99").Value = 794.71429
$ws.Range("N99").Value = -23742.2
$ws.Range("H101").Value = 300
$ws.Range("I101").Value = 450
$ws.Range("J101").Value = 150
$ws.Range("K101").Value = 1350
$ws.Range("L101").Value = 450
$ws.Range("M101").Value = 272
$ws.Range("N101").Value = -3694
$ws.Range("H102").Value = 99999.8
$ws.Range("J102").Value = 99999.8
$ws.Range("L102").Value = 99999.8
$ws.Range("N102").Value = -106489.8
$ws.Range("H105").Value = 24549.5
$ws.Range("J105").Value = 24549.5
$ws.Range("L105").Value = 24549.5
$ws.Range("N105").Value = -31537.5
$ws.Range("H111").Value = 1468.4286
$ws.Range("I111").Value = 1471.6
$ws.Range("K111").Value = 4414.799999999999
$ws.Range("M111").Value = -1347.799999999999
$ws.Range("H116").Value = 6307.5454
$ws.Range("I116").Value = 6377.2
$ws.Range("J116").Value = 6249.5
$ws.Range("K116").Value = 6377.2
$ws.Range("L116").Value = 6249.5
$ws.Range("M116").Value = -2935.2
$ws.Range("N116").Value = -13133.5
$ws.Range("H129").Value = 2001.6
$ws.Range("I129").Value = 2001.6
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 6004.799999999999
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -1004.799999999999
$ws.Range("N129").ClearContents()
$ws.Range("H138").Value = 5169.7637
$ws.Range("I138").Value = 3888.0527
$ws.Range("J138").Value = 5846.222
$ws.Range("K138").Value = 11664.1581
$ws.Range("L138").Value = 17538.666
$ws.Range("M138").Value = -6524.158100000001
$ws.Range("N138").Value = -27818.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26844.49
$ws.Range("I32").Value = 27524.5
$ws.Range("K32").Value = 27524.5
$ws.Range("M32").Value = -27237.5
$ws.Range("H45").Value = 3045.3333
$ws.Range("I45").Value = 1122.091
$ws.Range("J45").Value = 5160.9
$ws.Range("K45").Value = 1122.091
$ws.Range("L45").Value = 5160.9
$ws.Range("M45").Value = -745.0909999999999
$ws.Range("N45").Value = -5914.9
$ws.Range("H122").Value = 3666.6667
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H132").Value = 46888
$ws.Range("I132").Value = 55154.79
$ws.Range("K132").Value = 165464.37
$ws.Range("M132").Value = -162934.37

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2373
$ws.Range("I86").Value = 1807.1666
$ws.Range("J86").Value = 3221.75
$ws.Range("K86").Value = 1807.1666
$ws.Range("L86").Value = 3221.75
$ws.Range("M86").Value = -684.1666
$ws.Range("N86").Value = -5467.75
$ws.Range("H89").Value = 2373
$ws.Range("I89").Value = 1807.1666
$ws.Range("J89").Value = 3221.75
$ws.Range("K89").Value = 9035.833000000001
$ws.Range("L89").Value = 16108.75
$ws.Range("M89").Value = -3419.833000000001
$ws.Range("N89").Value = -27340.75
$ws.Range("H103").Value = 32643.25
$ws.Range("J103").Value = 32643.25
$ws.Range("L103").Value = 32643.25
$ws.Range("N103").Value = -34987.25
$ws.Range("H134").Value = 1884.697
$ws.Range("I134").Value = 1574.9231
$ws.Range("J134").Value = 3035.2856
$ws.Range("K134").Value = 4724.7693
$ws.Range("L134").Value = 9105.856800000001
$ws.Range("M134").Value = -2189.7693
$ws.Range("N134").Value = -14175.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H141").Value = 384250.12
$ws.Range("J141").Value = 494333.5
$ws.Range("L141").Value = 494333.5
$ws.Range("N141").Value = -504693.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 3461.182
$ws.Range("I99").Value = 1012.3333
$ws.Range("J99").Value = 6399.8
$ws.Range("K99").Value = 3036.9999
$ws.Range("L99").Value = 19199.4
$ws.Range("M99").Value = -790.9998999999998
$ws.Range("N99").Value = -23691.4
$ws.Range("H117").Value = 2739.6
$ws.Range("I117").Value = 1509.2858
$ws.Range("J117").Value = 3816.125
$ws.Range("K117").Value = 4527.857400000001
$ws.Range("L117").Value = 11448.375
$ws.Range("M117").Value = -1085.857400000001
$ws.Range("N117").Value = -18332.375
$ws.Range("H123").Value = 2140.8572
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H133").Value = 6505.222
$ws.Range("I133").Value = 2886.75
$ws.Range("J133").Value = 9400
$ws.Range("K133").Value = 8660.25
$ws.Range("L133").Value = 28200
$ws.Range("M133").Value = -3600.25
$ws.Range("N133").Value = -38320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4749.5713
$ws.Range("I70").Value = 4538.778
$ws.Range("J70").Value = 5129
$ws.Range("K70").Value = 4538.778
$ws.Range("L70").Value = 5129
$ws.Range("M70").Value = -4268.778
$ws.Range("N70").Value = -5669
$ws.Range("H73").Value = 4749.5713
$ws.Range("I73").Value = 4538.778
$ws.Range("J73").Value = 5129
$ws.Range("K73").Value = 4538.778
$ws.Range("L73").Value = 5129
$ws.Range("M73").Value = -3602.778
$ws.Range("N73").Value = -7001
$ws.Range("H107").Value = 72613.64
$ws.Range("I107").Value = 143349.14
$ws.Range("J107").Value = 1878.1428
$ws.Range("K107").Value = 143349.14
$ws.Range("L107").Value = 1878.1428
$ws.Range("M107").Value = -141429.14
$ws.Range("N107").Value = -5718.1428
$ws.Range("H126").Value = 5076.154
$ws.Range("I126").Value = 4237.1904
$ws.Range("K126").Value = 12711.5712
$ws.Range("M126").Value = -10241.5712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2102
$ws.Range("I40").Value = 2102
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2102
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1966
$ws.Range("N40").ClearContents()
$ws.Range("H122").Value = 4012.878
$ws.Range("I122").Value = 3417.9443
$ws.Range("J122").Value = 4478.478
$ws.Range("K122").Value = 10253.8329
$ws.Range("L122").Value = 13435.434
$ws.Range("M122").Value = -7803.832900000001
$ws.Range("N122").Value = -18335.434
$ws.Range("H132").Value = 50862.6
$ws.Range("I132").Value = 59636.43
$ws.Range("J132").Value = 4800
$ws.Range("K132").Value = 178909.29
$ws.Range("L132").Value = 14400
$ws.Range("M132").Value = -176379.29
$ws.Range("N132").Value = -19460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3475.1765
$ws.Range("I81").Value = 629.63635
$ws.Range("J81").Value = 8692
$ws.Range("K81").Value = 1259.2727
$ws.Range("L81").Value = 17384
$ws.Range("M81").Value = -198.2727
$ws.Range("N81").Value = -19506
$ws.Range("H84").Value = 3475.1765
$ws.Range("I84").Value = 629.63635
$ws.Range("J84").Value = 8692
$ws.Range("K84").Value = 6296.363499999999
$ws.Range("L84").Value = 86920
$ws.Range("M84").Value = -992.3634999999995
$ws.Range("N84").Value = -97528
$ws.Range("H107").Value = 1099.125
$ws.Range("I107").Value = 1133.3334
$ws.Range("J107").Value = 996.5
$ws.Range("K107").Value = 3400.0002
$ws.Range("L107").Value = 2989.5
$ws.Range("M107").Value = -1480.0002
$ws.Range("N107").Value = -6829.5
$ws.Range("H132").Value = 140944.06
$ws.Range("I132").Value = 152783.77
$ws.Range("K132").Value = 458351.3099999999
$ws.Range("M132").Value = -455821.3099999999

